$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 27 with formula in column B
$ws.Range("B27").Formula = "=60*1.6235+20.107"

# Move selection to B28 (the cell below the new data), matching post-edit state
$ws.Range("B28").Select()
